$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Keegan Murray"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Sacramento Kings"

$ws.Range("A8").Value = "Karl-Anthony Towns"
$ws.Range("B8").Value = "PF,C"
$ws.Range("C8").Value = "New York Knicks"

$ws.Range("A9").Value = "Peyton Watson"
$ws.Range("B9").Value = "SF,PF"
$ws.Range("C9").Value = "Denver Nuggets"

$ws.Range("A13").Value = "Trey Murphy III"
$ws.Range("B13").Value = "SF,PF"
$ws.Range("C13").Value = "New Orleans Pelicans"

$ws.Range("A14").Value = "OG Anunoby"
$ws.Range("B14").Value = "SF,PF"
$ws.Range("C14").Value = "New York Knicks"

$ws.Range("A15").Value = "Franz Wagner"
$ws.Range("B15").Value = "SF,PF"
$ws.Range("C15").Value = "Orlando Magic"

$ws.Range("A16").Value = "Jarrett Allen"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Cleveland Cavaliers"

$ws.Range("A19").Value = "Darius Garland"
$ws.Range("B19").Value = "PG"
$ws.Range("C19").Value = "Cleveland Cavaliers"
